$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new values. A leading apostrophe forces text entry (mirrors typing
# '0.9978 into Excel) so purely-numeric-looking strings like "0.9978" or
# "13.00" are NOT reinterpreted as numbers, matching the source file's
# inline-string cell type.
$ws.Range("D2").Value = '''23.534.28'
$ws.Range("E2").Value = '''  +2.22%  '
$ws.Range("D3").Value = '''1.636.69'
$ws.Range("E3").Value = '''  +3.72%  '
$ws.Range("D4").Value = '''0.9978'
$ws.Range("E4").Value = '''  -0.36%  '
$ws.Range("D5").Value = '''307.83'
$ws.Range("E5").Value = '''  +2.85%  '
$ws.Range("D6").Value = '''0.9971'
$ws.Range("E6").Value = '''  -0.52%  '
$ws.Range("D7").Value = '''0.3784'
$ws.Range("E7").Value = '''  +1.16%  '
$ws.Range("D8").Value = '''53.22'
$ws.Range("E8").Value = '''  +6.70%  '
$ws.Range("D9").Value = '''0.3684'
$ws.Range("E9").Value = '''  +3.50%  '
$ws.Range("D10").Value = '''1.287'
$ws.Range("E10").Value = '''  +6.14%  '
$ws.Range("D11").Value = '''0.08200'
$ws.Range("E11").Value = '''  +2.96%  '
$ws.Range("D12").Value = '''0.9980'
$ws.Range("E12").Value = '''  -0.32%  '
$ws.Range("D13").Value = '''23.34'
$ws.Range("E13").Value = '''  +7.14%  '
$ws.Range("D14").Value = '''6.676'
$ws.Range("E14").Value = '''  +3.85%  '
$ws.Range("D15").Value = '''0.00001285'
$ws.Range("E15").Value = '''  +5.57%  '
$ws.Range("D16").Value = '''7.484'
$ws.Range("E16").Value = '''  +2.88%  '
$ws.Range("D17").Value = '''1.631.58'
$ws.Range("E17").Value = '''  +3.08%  '
$ws.Range("D18").Value = '''95.01'
$ws.Range("E18").Value = '''  +3.62%  '
$ws.Range("D19").Value = '''0.06952'
$ws.Range("E19").Value = '''  +3.15%  '
$ws.Range("D20").Value = '''18.47'
$ws.Range("E20").Value = '''  +4.32%  '
$ws.Range("D21").Value = '''6.607'
$ws.Range("E21").Value = '''  +3.97%  '
$ws.Range("D22").Value = '''0.9971'
$ws.Range("E22").Value = '''  -0.46%  '
$ws.Range("D23").Value = '''23.543.73'
$ws.Range("E23").Value = '''  +2.33%  '
$ws.Range("D24").Value = '''13.00'
$ws.Range("E24").Value = '''  +2.96%  '
$ws.Range("D25").Value = '''3.145'
$ws.Range("E25").Value = '''  +11.85%  '
$ws.Range("D26").Value = '''2.421'
$ws.Range("E26").Value = '''  +2.48%  '
$ws.Range("E27").Value = '''  +4.31%  '
$ws.Range("D28").Value = '''151.53'
$ws.Range("D29").Value = '''5.315'
$ws.Range("E29").Value = '''  +2.65%  '
$ws.Range("D30").Value = '''137.11'
$ws.Range("E30").Value = '''  +4.40%  '
$ws.Range("D31").Value = '''2.429'
$ws.Range("E31").Value = '''  +4.52%  '
$ws.Range("D32").Value = '''6.879'
$ws.Range("E32").Value = '''  +6.35%  '
$ws.Range("D33").Value = '''1.810.86'
$ws.Range("E33").Value = '''  +3.09%  '
$ws.Range("D34").Value = '''0.9780'
$ws.Range("E34").Value = '''  +5.47%  '
$ws.Range("D35").Value = '''0.02832'
$ws.Range("E35").Value = '''  +6.70%  '
$ws.Range("D36").Value = '''10.45'
$ws.Range("E36").Value = '''  +5.63%  '
$ws.Range("D37").Value = '''0.07495'
$ws.Range("E37").Value = '''  +2.37%  '
$ws.Range("D38").Value = '''6.238'
$ws.Range("E38").Value = '''  +4.83%  '
$ws.Range("D39").Value = '''0.2543'
$ws.Range("E39").Value = '''  +2.77%  '
$ws.Range("D40").Value = '''0.08873'
$ws.Range("E40").Value = '''  +1.55%  '
$ws.Range("D41").Value = '''1.399'
$ws.Range("E41").Value = '''  +4.54%  '
$ws.Range("D42").Value = '''0.7177'
$ws.Range("E42").Value = '''  +5.19%  '
$ws.Range("D43").Value = '''12.75'
$ws.Range("E43").Value = '''  +8.69%  '
$ws.Range("D44").Value = '''16.22'
$ws.Range("E44").Value = '''  +10.69%  '
$ws.Range("D45").Value = '''0.6655'
$ws.Range("E45").Value = '''  +5.36%  '
$ws.Range("D46").Value = '''2.367'
$ws.Range("E46").Value = '''  +5.98%  '
$ws.Range("D47").Value = '''4.037'
$ws.Range("E47").Value = '''  +1.82%  '
$ws.Range("D48").Value = '''0.9963'
$ws.Range("E48").Value = '''  -0.46%  '
$ws.Range("D49").Value = '''0.08057'
$ws.Range("E49").Value = '''  +2.73%  '
$ws.Range("D50").Value = '''132.14'
$ws.Range("E50").Value = '''  +0.91%  '
$ws.Range("D51").Value = '''1.216'
$ws.Range("E51").Value = '''  +3.15%  '

# Strip the quote-prefix formatting flag Excel sets when a value is entered
# with a leading apostrophe, so the cells' style stays identical to the
# untouched cells around them (no explicit style index).
$ws.Range("D2:E51").ClearFormats()

